# Add a new "chemical_recycling_pyrolysis" parameter row right after the
# existing "chemical_recycling_gasification" row (row 9), pushing every
# row below it down by one. This mirrors the commit:
#   "revision, added pyrolysis and additional figures"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 10 (shifts old rows 10..24 to 11..25,
# dimension grows from A1:C24 to A1:C25 automatically).
$ws.Rows("10").Insert()

# Populate the new row: parameter name + boolean value (TRUE),
# matching the layout of the other boolean-flag rows (e.g. row 9).
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
